$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# ---------------------------------------------------------------------------
# Rows 2-6 (years 2014/12 .. 2018/12): overwrite the financial figures with
# the corrected values. Columns J and O (only present on rows 2-5) are
# cleared out entirely since they no longer exist after the fix.
# ---------------------------------------------------------------------------

# Row 2 - 2014/12
$ws.Range("D2").Value  = 722
$ws.Range("E2").Value  = -55
$ws.Range("F2").Value  = -55
$ws.Range("G2").Value  = -62
$ws.Range("H2").Value  = -62
$ws.Range("I2").Value  = -62
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value  = 618
$ws.Range("L2").Value  = 179
$ws.Range("M2").Value  = 439
$ws.Range("N2").Value  = 439
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value  = 138
$ws.Range("Q2").Value  = -31
$ws.Range("R2").Value  = 22
$ws.Range("S2").Value  = -1
$ws.Range("T2").Value  = 4
$ws.Range("U2").Value  = -35
$ws.Range("V2").Value  = 1
$ws.Range("W2").Value  = -7.58
$ws.Range("X2").Value  = -8.609999999999999
$ws.Range("Y2").Value  = -13.21
$ws.Range("Z2").Value  = -9.220000000000001
$ws.Range("AA2").Value = 40.66
$ws.Range("AB2").Value = 350.04
$ws.Range("AC2").Value = -1124
$ws.Range("AD2").Value = -4.56
$ws.Range("AE2").Value = 9977
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 5530735

# Row 3 - 2015/12
$ws.Range("D3").Value  = 545
$ws.Range("E3").Value  = -59
$ws.Range("F3").Value  = -59
$ws.Range("G3").Value  = -58
$ws.Range("H3").Value  = -58
$ws.Range("I3").Value  = -58
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value  = 548
$ws.Range("L3").Value  = 161
$ws.Range("M3").Value  = 387
$ws.Range("N3").Value  = 387
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value  = 138
$ws.Range("Q3").Value  = -42
$ws.Range("R3").Value  = -3
$ws.Range("S3").Value  = 29
$ws.Range("T3").Value  = 7
$ws.Range("U3").Value  = -49
$ws.Range("V3").Value  = 30
$ws.Range("W3").Value  = -10.75
$ws.Range("X3").Value  = -10.7
$ws.Range("Y3").Value  = -14.13
$ws.Range("Z3").Value  = -10.01
$ws.Range("AA3").Value = 41.74
$ws.Range("AB3").Value = 307.83
$ws.Range("AC3").Value = -1055
$ws.Range("AD3").Value = -5.83
$ws.Range("AE3").Value = 8786
$ws.Range("AF3").Value = 0.7
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5530735

# Row 4 - 2016/12
$ws.Range("D4").Value  = 411
$ws.Range("E4").Value  = -80
$ws.Range("F4").Value  = -80
$ws.Range("G4").Value  = 30
$ws.Range("H4").Value  = 30
$ws.Range("I4").Value  = 30
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value  = 639
$ws.Range("L4").Value  = 197
$ws.Range("M4").Value  = 442
$ws.Range("N4").Value  = 442
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value  = 138
$ws.Range("Q4").Value  = -89
$ws.Range("R4").Value  = 79
$ws.Range("S4").Value  = 18
$ws.Range("T4").Value  = 6
$ws.Range("U4").Value  = -95
$ws.Range("V4").Value  = 35
$ws.Range("W4").Value  = -19.47
$ws.Range("X4").Value  = 7.32
$ws.Range("Y4").Value  = 7.27
$ws.Range("Z4").Value  = 5.08
$ws.Range("AA4").Value = 44.6
$ws.Range("AB4").Value = 329.62
$ws.Range("AC4").Value = 545
$ws.Range("AD4").Value = 9.67
$ws.Range("AE4").Value = 9421
$ws.Range("AF4").Value = 0.5600000000000001
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 5530735

# Row 5 - 2017/12
$ws.Range("D5").Value  = 376
$ws.Range("E5").Value  = -114
$ws.Range("F5").Value  = -114
$ws.Range("G5").Value  = -114
$ws.Range("H5").Value  = -114
$ws.Range("I5").Value  = -114
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value  = 579
$ws.Range("L5").Value  = 252
$ws.Range("M5").Value  = 327
$ws.Range("N5").Value  = 327
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value  = 138
$ws.Range("Q5").Value  = -110
$ws.Range("R5").Value  = 51
$ws.Range("S5").Value  = 60
$ws.Range("T5").Value  = 13
$ws.Range("U5").Value  = -123
$ws.Range("V5").Value  = 94
$ws.Range("W5").Value  = -30.37
$ws.Range("X5").Value  = -30.26
$ws.Range("Y5").Value  = -29.59
$ws.Range("Z5").Value  = -18.68
$ws.Range("AA5").Value = 77.12
$ws.Range("AB5").Value = 247.39
$ws.Range("AC5").Value = -2056
$ws.Range("AD5").Value = -2.02
$ws.Range("AE5").Value = 6949
$ws.Range("AF5").Value = 0.6
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 5530735

# Row 6 - 2018/12 (no J/O columns existed here already)
$ws.Range("D6").Value  = 408
$ws.Range("E6").Value  = -10
$ws.Range("F6").Value  = -10
$ws.Range("G6").Value  = 7
$ws.Range("H6").Value  = 7
$ws.Range("I6").Value  = 7
$ws.Range("K6").Value  = 641
$ws.Range("L6").Value  = 298
$ws.Range("M6").Value  = 343
$ws.Range("N6").Value  = 343
$ws.Range("P6").Value  = 138
$ws.Range("Q6").Value  = -79
$ws.Range("R6").Value  = -4
$ws.Range("S6").Value  = 67
$ws.Range("T6").Value  = 3
$ws.Range("U6").Value  = -82
$ws.Range("V6").Value  = 155
$ws.Range("W6").Value  = -2.4
$ws.Range("X6").Value  = 1.76
$ws.Range("Y6").Value  = 2.14
$ws.Range("Z6").Value  = 1.18
$ws.Range("AA6").Value = 86.76000000000001
$ws.Range("AB6").Value = 254.02
$ws.Range("AC6").Value = 130
$ws.Range("AD6").Value = 27.92
$ws.Range("AE6").Value = 7061
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 5530735

# ---------------------------------------------------------------------------
# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)): all the estimate figures
# were bad data and are removed entirely, leaving only the no./label/period
# columns (A, B, C).
# ---------------------------------------------------------------------------
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
